$d = $word.ActiveDocument

# 1. Header: チャイの売上合計 (単位) -> チャイの売上合計 (ユニット数)
$d.Content.Find.Execute("チャイの売上合計 (単位)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "チャイの売上合計 (ユニット数)", 1)

# 2. Header: 職人チャイ販売 (ユニット) -> Artisanal Chai の売上 (ユニット数)
$d.Content.Find.Execute("職人チャイ販売 (ユニット)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Artisanal Chai の売上 (ユニット数)", 1)

# 3. Header: 事前に作成されたチャイの売上 (単位) -> 既製チャイの売上 (ユニット数)
$d.Content.Find.Execute("事前に作成されたチャイの売上 (単位)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "既製チャイの売上 (ユニット数)", 1)

# 4 & 5. Header: ソーシャル メディア エンゲージメント -> ソーシャル メディア エンゲージメント (ビュー), and bold it
$d.Content.Find.Execute("ソーシャル メディア エンゲージメント", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ソーシャル メディア エンゲージメント (ビュー)", 1)

# Now apply bold formatting to that run: search again and set Bold
$rng = $d.Content
$rng.Find.Execute("ソーシャル メディア エンゲージメント (ビュー)", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
$rng.Font.Bold = $true

# 6. Date 3/31/2023 (2nd occurrence) -> 2023/5/31
$rng2 = $d.Content
$found1 = $rng2.Find.Execute("3/31/2023", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
# move past first occurrence, then replace the next one in the remaining range
$searchStart = $rng2.End
$rng3 = $d.Range($searchStart, $d.Content.End)
$rng3.Find.Execute("3/31/2023", $true, $false, $false, $false, $false,
                    $true, 1, $false, "2023/5/31", 1)

# 7. Date 3/30/2023 -> 2023/7/30
$d.Content.Find.Execute("3/30/2023", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023/7/30", 1)

# 8. Date 2023 年 9 月 5 日 -> 2023/9/30
$d.Content.Find.Execute("2023 年 9 月 5 日", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023/9/30", 1)

# 9. Date 2020/11/30 -> 2023/11/30
$d.Content.Find.Execute("2020/11/30", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023/11/30", 1)
